# Auto-generated edit script: updates crypto price/volume values and swaps Filecoin/ImmutableX rows 31-32
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '36.336.97'
$ws.Range('E2').Value = '  -1.51%  '
$ws.Range('D3').Value = '2.035.41'
$ws.Range('E3').Value = '  -0.56%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '244.44'
$ws.Range('E5').Value = '  -0.51%  '
$ws.Range('D6').Value = '0.654'
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '''52.70'
$ws.Range('E8').Value = '  -8.41%  '
$ws.Range('D9').Value = '60.74'
$ws.Range('E9').Value = '  +3.05%  '
$ws.Range('D10').Value = '0.357'
$ws.Range('E10').Value = '  -3.54%  '
$ws.Range('D11').Value = '0.0734'
$ws.Range('E11').Value = '  -5.52%  '
$ws.Range('E12').Value = '  -4.11%  '
$ws.Range('D13').Value = '0.917'
$ws.Range('E13').Value = '  +5.38%  '
$ws.Range('D14').Value = '14.33'
$ws.Range('E14').Value = '  -5.58%  '
$ws.Range('D15').Value = '2.337.90'
$ws.Range('E15').Value = '  -0.40%  '
$ws.Range('D16').Value = '5.32'
$ws.Range('E16').Value = '  -5.05%  '
$ws.Range('D17').Value = '2.053.18'
$ws.Range('E17').Value = '  +0.98%  '
$ws.Range('D18').Value = '36.271.17'
$ws.Range('E18').Value = '  -1.54%  '
$ws.Range('D19').Value = '16.68'
$ws.Range('E19').Value = '  -6.59%  '
$ws.Range('D20').Value = '70.75'
$ws.Range('E20').Value = '  -3.57%  '
$ws.Range('D21').Value = '0.0₃0838'
$ws.Range('E21').Value = '  -5.27%  '
$ws.Range('D22').Value = '234.76'
$ws.Range('E22').Value = '  -0.32%  '
$ws.Range('D23').Value = '5.09'
$ws.Range('E23').Value = '  -5.34%  '
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').Value = '2.35'
$ws.Range('E25').Value = '  -3.84%  '
$ws.Range('D26').Value = '2.21'
$ws.Range('E26').Value = '  +1.01%  '
$ws.Range('D27').Value = '162.85'
$ws.Range('E27').Value = '  -3.16%  '
$ws.Range('E28').Value = '  -12.23%  '
$ws.Range('D29').Value = '19.62'
$ws.Range('E29').Value = '  -1.44%  '
$ws.Range('E30').Value = '  -2.87%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '1.16'
$ws.Range('E31').Value = '  +5.54%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '4.88'
$ws.Range('E32').Value = '  -10.21%  '
$ws.Range('D33').Value = '0.0582'
$ws.Range('E33').Value = '  -4.44%  '
$ws.Range('D34').Value = '4.32'
$ws.Range('E34').Value = '  -7.14%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').Value = '0.0861'
$ws.Range('E36').Value = '  +4.95%  '
$ws.Range('E37').Value = '  -0.77%  '
$ws.Range('D38').Value = '2.16'
$ws.Range('E38').Value = '  -7.09%  '
$ws.Range('D39').Value = '4.89'
$ws.Range('E39').Value = '  -4.73%  '
$ws.Range('D40').Value = '''1.20'
$ws.Range('E40').Value = '  -8.23%  '
$ws.Range('E41').Value = '  -4.36%  '
$ws.Range('D42').Value = '''0.0210'
$ws.Range('E42').Value = '  -5.49%  '
$ws.Range('E43').Value = '  -5.11%  '
$ws.Range('D44').Value = '91.69'
$ws.Range('E44').Value = '  -4.80%  '
$ws.Range('D45').Value = '0.0887'
$ws.Range('E45').Value = '  -5.74%  '
$ws.Range('D46').Value = '1.367.86'
$ws.Range('E46').Value = '  +5.54%  '
$ws.Range('D47').Value = '7.31'
$ws.Range('E47').Value = '  +8.79%  '
$ws.Range('D48').Value = '15.38'
$ws.Range('E48').Value = '  -8.59%  '
$ws.Range('D49').Value = '''2.90'
$ws.Range('E49').Value = '  +1.68%  '
$ws.Range('D50').Value = '2.219.16'
$ws.Range('E50').Value = '  -0.66%  '
$ws.Range('E51').Value = '  -5.33%  '

$wb.Save()
